# Update "想去人数" (want-to-go count) values in column F
# for the "展览" (sheet 1) and "全部类型" (sheet 4) worksheets,
# reflecting refreshed counts fetched at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 1163
$ws1.Range("F4").Value  = 143
$ws1.Range("F6").Value  = 281
$ws1.Range("F7").Value  = 96
$ws1.Range("F8").Value  = 1215
$ws1.Range("F9").Value  = 17220
$ws1.Range("F10").Value = 313
$ws1.Range("F13").Value = 6575
$ws1.Range("F15").Value = 143
$ws1.Range("F16").Value = 91
$ws1.Range("F17").Value = 45
$ws1.Range("F19").Value = 1286
$ws1.Range("F20").Value = 135
$ws1.Range("F28").Value = 88
$ws1.Range("F32").Value = 11676
$ws1.Range("F33").Value = 1262
$ws1.Range("F36").Value = 246
$ws1.Range("F39").Value = 83

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 1163
$ws4.Range("F4").Value  = 143
$ws4.Range("F6").Value  = 281
$ws4.Range("F7").Value  = 96
$ws4.Range("F8").Value  = 1215
$ws4.Range("F9").Value  = 17220
$ws4.Range("F10").Value = 313
$ws4.Range("F13").Value = 6575
$ws4.Range("F15").Value = 143
$ws4.Range("F16").Value = 91
$ws4.Range("F17").Value = 45
$ws4.Range("F19").Value = 1286
$ws4.Range("F20").Value = 135
$ws4.Range("F28").Value = 88
$ws4.Range("F33").Value = 11676
$ws4.Range("F34").Value = 1262
$ws4.Range("F37").Value = 246
$ws4.Range("F40").Value = 83

$wb.Save()
